$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) contains numeric-looking text (e.g. "28.125.22", "0.9974")
# that must remain literal text rather than be auto-converted to numbers by Excel.
# Force the whole Price column to Text format before writing the new values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.125.22'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '1.760.67'
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '335.02'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").Value = '0.9974'
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").Value = '0.3782'
$ws.Range("E7").Value = '  -3.28%  '
$ws.Range("D8").Value = '0.3362'
$ws.Range("E8").Value = '  -3.45%  '
$ws.Range("D9").Value = '45.68'
$ws.Range("E9").Value = '  -5.61%  '
$ws.Range("D10").Value = '1.126'
$ws.Range("E10").Value = '  -5.53%  '
$ws.Range("D11").Value = '0.07215'
$ws.Range("E11").Value = '  -4.70%  '
$ws.Range("D12").Value = '22.61'
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("D13").Value = '0.9996'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").Value = '6.193'
$ws.Range("E14").Value = '  -4.91%  '
$ws.Range("D15").Value = '7.215'
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("D16").Value = '1.758.15'
$ws.Range("E16").Value = '  -3.05%  '
$ws.Range("D17").Value = '0.00001055'
$ws.Range("E17").Value = '  -4.37%  '
$ws.Range("D18").Value = '0.06580'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").Value = '80.76'
$ws.Range("E19").Value = '  -4.97%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '17.01'
$ws.Range("E21").Value = '  -4.30%  '
$ws.Range("D22").Value = '6.277'
$ws.Range("E22").Value = '  -4.36%  '
$ws.Range("D23").Value = '28.118.14'
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D24").Value = '11.66'
$ws.Range("E24").Value = '  -6.46%  '
$ws.Range("D25").Value = '2.389'
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").Value = '153.54'
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("E27").Value = '  -6.43%  '
$ws.Range("D28").Value = '2.337'
$ws.Range("E28").Value = '  -7.60%  '
$ws.Range("D29").Value = '1.959.14'
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("D30").Value = '1.273'
$ws.Range("E30").Value = '  -14.61%  '
$ws.Range("D31").Value = '131.90'
$ws.Range("E31").Value = '  -2.69%  '
$ws.Range("D32").Value = '4.012'
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").Value = '5.820'
$ws.Range("E33").Value = '  -5.41%  '
$ws.Range("D34").Value = '0.08802'
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("D35").Value = '12.28'
$ws.Range("E35").Value = '  -5.86%  '
$ws.Range("D36").Value = '0.02343'
$ws.Range("E36").Value = '  -3.35%  '
$ws.Range("D37").Value = '0.6647'
$ws.Range("E37").Value = '  -4.50%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.06208'
$ws.Range("E38").Value = '  -5.07%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.175'
$ws.Range("E39").Value = '  -5.39%  '
$ws.Range("D40").Value = '0.2119'
$ws.Range("E40").Value = '  -4.30%  '
$ws.Range("D41").Value = '1.218'
$ws.Range("D42").Value = '1.456'
$ws.Range("E42").Value = '  -9.55%  '
$ws.Range("D43").Value = '8.032'
$ws.Range("E43").Value = '  -5.33%  '
$ws.Range("D44").Value = '0.9972'
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '13.72'
$ws.Range("E45").Value = '  -5.39%  '
$ws.Range("D46").Value = '0.6053'
$ws.Range("E46").Value = '  -5.99%  '
$ws.Range("D47").Value = '3.814'
$ws.Range("E47").Value = '  -1.46%  '
$ws.Range("D48").Value = '129.78'
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("D49").Value = '2.017'
$ws.Range("E49").Value = '  -6.29%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.188'
$ws.Range("E50").Value = '  +3.17%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.07211'
$ws.Range("E51").Value = '  +0.18%  '
